# Update month headers (row 1, columns C:N) from 2024 to 2025 on every sheet.
$wb = $excel.ActiveWorkbook

$months = @("Jan", "Feb", "Mar", "Apr", "May", "Jun", "Jul", "Aug", "Sep", "Oct", "Nov", "Dec")

foreach ($ws in $wb.Worksheets) {
    for ($i = 0; $i -lt $months.Length; $i++) {
        $col = 3 + $i   # Column C is index 3
        $cell = $ws.Cells.Item(1, $col)
        $cell.Value = $months[$i] + "25"
    }
}
